# weekly report generation done
# Add a new row to Sheet1 for the "weeklyRepCreationServiceUrl" config entry,
# mirroring the existing rows (label in column A, hyperlinked "http://google.com"
# in column B using the sheet's existing Hyperlink cell style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 16: label + hyperlinked URL (same pattern as rows 4-15)
$ws.Range("A16").Value = "weeklyRepCreationServiceUrl"
$ws.Range("B16").Value = "http://google.com"

$null = $ws.Hyperlinks.Add($ws.Range("B16"), "http://google.com/")
$ws.Range("B16").Style = "Hyperlink"

# Move the selection to A17, matching where the cursor ends up after filling row 16
$null = $ws.Range("A17").Select()
